# Fruta / hortaliza, semanal
# The weekly data refresh reshuffles the daily price rows (2-14) of the
# "Alcachofa" sheet: each row keeps its fixed identity columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Clasificación) but receives a different set of values for the
# date/quality/volume/price/origin columns (D, I, J, K, L, M, N, O, P, Q).
#
# Capture the current ("before") values for those columns on every data
# row, then redistribute them onto the rows according to the mapping
# observed between the previous and the new workbook revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current values of the columns that get reshuffled.
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $rowValues[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowValues
}

# Maps each destination row to the source row whose values it should
# now hold (derived from comparing the workbook before/after the edit).
$mapping = @{
    2  = 6
    3  = 2
    4  = 10
    5  = 13
    6  = 9
    7  = 3
    8  = 11
    9  = 12
    10 = 8
    11 = 7
    12 = 14
    13 = 4
    14 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcValues[$c]
    }
}
